# fix xlsx2csv component, manage empty cell values
#
# Populate the third sheet ("Sheet1", sheet3.xml) with a small sparse
# table that exercises empty/missing cells in the middle of rows -
# exactly the xlsx2csv scenario the commit message refers to.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sparse data - note the gaps (empty cells) on each row/column on purpose.
$ws.Range("A1").Value = "a"
$ws.Range("C1").Value = "c"

$ws.Range("A2").Value = "a"
$ws.Range("B2").Value = "b"

$ws.Range("A3").Value = "a"
$ws.Range("D3").Value = "d"

$ws.Range("E4").Value = "e"

# Outline levels (sheetFormatPr outlineLevelRow / outlineLevelCol).
$ws.Rows(1).OutlineLevel = 3
$ws.Columns(1).OutlineLevel = 4

# Move/restore the on-screen selection.
[void]$ws.Range("G9").Select()
